$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 324
$ws.Range("G3").Value = 36
$ws.Range("G4").Value = 332
$ws.Range("G5").Value = 200
$ws.Range("G6").Value = 1556
$ws.Range("G7").Value = 7
$ws.Range("G9").Value = 789
$ws.Range("G10").Value = 66
$ws.Range("G12").Value = 403
